$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "58.926.36"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.307.81"
$ws.Range("E3").Value = "  -4.73%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "550.46"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "131.32"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").Value = "2.304.66"
$ws.Range("E9").Value = "  -4.81%  "
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").Value = "24.06"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "2.714.86"
$ws.Range("E15").Value = "  -4.93%  "
$ws.Range("D16").Value = "58.895.75"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "2.217.88"
$ws.Range("E18").Value = "  -8.50%  "
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").Value = "4.34"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("D21").Value = "316.20"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "63.22"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "8.13"
$ws.Range("E27").Value = "  -6.18%  "
$ws.Range("E28").Value = "  -7.70%  "
$ws.Range("D29").Value = "1.77"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "169.94"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "0.0₃0734"
$ws.Range("E31").Value = "  -5.32%  "
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("D33").Value = "5.84"
$ws.Range("E33").Value = "  -4.60%  "
$ws.Range("D34").Value = "0.385"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "17.84"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -6.10%  "
$ws.Range("E39").Value = "  -5.56%  "
$ws.Range("D40").Value = "38.08"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").Value = "306.95"
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("D43").Value = "141.57"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").Value = "3.47"
$ws.Range("E44").Value = "  -5.21%  "
$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "18.86"
$ws.Range("E47").Value = "  -4.67%  "
$ws.Range("D48").Value = "0.560"
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("D49").Value = "0.0217"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "16.77"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("E51").Value = "  -0.15%  "
